$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Month number -> abbreviated Spanish month name (as used by the workbook author)
$monthNames = @{
    1  = "Ene."
    2  = "Feb."
    3  = "Mar."
    4  = "Abr."
    5  = "May."
    6  = "Jun."
    7  = "Jul."
    8  = "Ago."
    9  = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

# Column C ("Mes") holds numeric month values 1-12 for rows 6 through 80;
# replace each with its abbreviated month-name text equivalent.
for ($r = 6; $r -le 80; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $monthNum = $cell.Value2
    if ($monthNum -ne $null -and $monthNames.ContainsKey([int]$monthNum)) {
        $cell.Value = $monthNames[[int]$monthNum]
    }
}

# Move the window position to match the saved workbook view (xWindow changed
# from -120 to 23880).
$excel.ActiveWindow.Left = 23880
